$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation for column D (and B/C for swapped rows) so that
# numeric-looking values (e.g. "0.9989", "1.000") are stored as text, matching
# the original inline-string cell types.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.505.59"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.687.42"
$ws.Range("E3").Value = "  +2.51%  "
$ws.Range("D4").Value = "0.9989"
$ws.Range("E4").Value = "  -0.49%  "
$ws.Range("D5").Value = "218.01"
$ws.Range("E5").Value = "  +5.20%  "
$ws.Range("D6").Value = "0.5343"
$ws.Range("E6").Value = "  +3.40%  "
$ws.Range("D7").Value = "0.9998"
$ws.Range("E7").Value = "  -0.41%  "
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  +3.84%  "
$ws.Range("D9").Value = "0.06460"
$ws.Range("E9").Value = "  +3.61%  "
$ws.Range("D10").Value = "21.48"
$ws.Range("E10").Value = "  +4.08%  "
$ws.Range("D11").Value = "0.07793"
$ws.Range("E11").Value = "  +2.97%  "
$ws.Range("D12").Value = "1.700.73"
$ws.Range("E12").Value = "  +3.29%  "
$ws.Range("D13").Value = "4.516"
$ws.Range("E13").Value = "  +3.31%  "
$ws.Range("D14").Value = "0.5669"
$ws.Range("E14").Value = "  +5.77%  "
$ws.Range("D15").Value = ("0.0{0}8457" -f [char]0x2085)
$ws.Range("E15").Value = "  +6.06%  "
$ws.Range("D16").Value = "66.37"
$ws.Range("E16").Value = "  +0.69%  "
$ws.Range("D17").Value = "26.498.57"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("D18").Value = "4.845"
$ws.Range("E18").Value = "  +4.02%  "
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.38%  "
$ws.Range("D20").Value = "195.31"
$ws.Range("E20").Value = "  +4.63%  "
$ws.Range("D21").Value = "10.43"
$ws.Range("E21").Value = "  +4.01%  "
$ws.Range("D22").Value = "6.425"
$ws.Range("E22").Value = "  +5.12%  "
$ws.Range("D23").Value = "1.000"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("D24").Value = "143.64"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").Value = "0.1275"
$ws.Range("E25").Value = "  +6.00%  "
$ws.Range("D26").Value = "7.517"
$ws.Range("E26").Value = "  +2.30%  "
$ws.Range("D27").Value = "16.31"
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("D28").Value = "1.432"
$ws.Range("E28").Value = "  +4.34%  "
$ws.Range("D29").Value = "0.06198"
$ws.Range("E29").Value = "  +3.08%  "
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  +2.77%  "
$ws.Range("D31").Value = "3.572"
$ws.Range("E31").Value = "  +4.31%  "
$ws.Range("D32").Value = "3.469"
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("D33").Value = "1.724"
$ws.Range("E33").Value = "  +6.61%  "
$ws.Range("D34").Value = "1.025"
$ws.Range("E34").Value = "  +5.57%  "
$ws.Range("D35").Value = "2.804"
$ws.Range("E35").Value = "  +2.53%  "
$ws.Range("D36").Value = "2.405"
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").Value = "0.5776"
$ws.Range("E37").Value = "  -1.52%  "
$ws.Range("D38").Value = "0.01652"
$ws.Range("E38").Value = "  +3.66%  "
$ws.Range("D39").Value = "5.970"
$ws.Range("E39").Value = "  +2.16%  "
$ws.Range("D40").Value = "0.8703"
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("D41").Value = "1.059.78"
$ws.Range("E41").Value = "  -1.16%  "
$ws.Range("D42").Value = "1.000"
$ws.Range("E42").Value = "  -0.33%  "
$ws.Range("D43").Value = "100.39"
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "1.835.55"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("D45").Value = "57.36"
$ws.Range("E45").Value = "  +5.29%  "
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").Value = "8.174"
$ws.Range("E47").Value = "  +2.36%  "
$ws.Range("D48").Value = "1.003"
$ws.Range("E48").Value = "  +0.61%  "

# Rows 49-51: Cronos/Aptos swap places (with updated price/volume), and Mantle
# is replaced by Algorand.
$ws.Range("B49").Value = "Aptos"
$ws.Range("C49").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D49").Value = "6.106"
$ws.Range("E49").Value = "  +4.53%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").Value = "0.05211"
$ws.Range("E50").Value = "  -0.04%  "
$ws.Range("B51").Value = "Algorand"
$ws.Range("C51").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D51").Value = "0.09959"
$ws.Range("E51").Value = "  +3.49%  "

# Reset style back to Normal so we don't leave a stray "@" number-format
# applied to the cells (value already committed as text above).
$ws.Range("D2:D51").Style = "Normal"
